# Rename the unclear "dependents" month-abbreviation headers (JAN..DEC)
# on the "Pool" sheet to dependents_1 .. dependents_12 (fixes #13)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool")

$labels = @(
    "dependents_1", "dependents_2", "dependents_3", "dependents_4",
    "dependents_5", "dependents_6", "dependents_7", "dependents_8",
    "dependents_9", "dependents_10", "dependents_11", "dependents_12"
)

# These headers live in row 1, columns N through Y (14..25)
$startCol = 14
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $labels[$i]
}
